$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 265, shifting existing rows (265-283) down to (266-284),
# carrying their formatting (including the date style on column D) with them.
$ws.Rows("265:265").Insert()

# Populate the newly inserted row 265 with the new record's values.
$ws.Range("A265").Value = 4
$ws.Range("B265").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C265").Value = "Los Lagos"
$ws.Range("D265").Value = 44746
$ws.Range("D265").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E265").Value = 10
$ws.Range("F265").Value = 100112043
$ws.Range("G265").Value = "Pepino ensalada"
$ws.Range("H265").Value = "Sin especificar"
$ws.Range("I265").Value = "Primera"
$ws.Range("J265").Value = 150
$ws.Range("K265").Value = 22000
$ws.Range("L265").Value = 23000
$ws.Range("M265").Value = 22467
$ws.Range("N265").Value = "$/caja 60 unidades"
$ws.Range("O265").Value = "Región de Arica y Parinacota"
$ws.Range("P265").Value = 374
$ws.Range("Q265").Value = 60
$ws.Range("R265").Value = "Hortaliza"
